# Add a new error code row to the SQS error codes worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 14: Code / Description
$ws.Cells.Item(14, 1).Value = 412112
$ws.Cells.Item(14, 2).Value = "The ARN for the specified SQS queue could not be found."

# Match formatting of the existing data rows (copy style from row 13, the
# last existing data row, onto the newly added row 14), scoped to just the
# two used columns so we don't spray styles across the whole sheet.
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Re-set the values since PasteSpecial(formats) should not have touched them,
# but make sure they are exactly right regardless.
$ws.Cells.Item(14, 1).Value = 412112
$ws.Cells.Item(14, 2).Value = "The ARN for the specified SQS queue could not be found."

# Update the active selection to match the end-state recorded in the diff.
$ws.Range("B15").Select()
